# Updates Ligand/Receptor average & total expression values and derived
# specificity / edge-weight metrics in rows 2-10 to reflect the new TPM
# based recomputation (see commit message "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"1.046867666666667"
$ws.Range("H2").Value = [double]"3.140603"
$ws.Range("I2").Value = [double]"0.000687505225377314"
$ws.Range("J2").Value = [double]"0.000687505225377314"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.2615913333333333"
$ws.Range("N2").Value = [double]"0.784774"
$ws.Range("O2").Value = [double]"0.08239613548481725"
$ws.Range("P2").Value = [double]"0.08239613548481727"
$ws.Range("Q2").Value = [double]"0.2738515087468889"
$ws.Range("R2").Value = [double]"2.464663578722"
$ws.Range("S2").Value = [double]"5.664777369670898E-05"
$ws.Range("T2").Value = [double]"5.664777369670899E-05"

# Row 3
$ws.Range("G3").Value = [double]"1.046867666666667"
$ws.Range("H3").Value = [double]"3.140603"
$ws.Range("I3").Value = [double]"0.000687505225377314"
$ws.Range("J3").Value = [double]"0.000687505225377314"
$ws.Range("N3").Value = [double]"5.233242000000001"
$ws.Range("O3").Value = [double]"0.5494561706387266"
$ws.Range("P3").Value = [double]"0.5494561706387268"
$ws.Range("Q3").Value = [double]"1.826170613880667"
$ws.Range("R3").Value = [double]"16.435535524926"
$ws.Range("S3").Value = [double]"0.0003777539884299336"
$ws.Range("T3").Value = [double]"0.0003777539884299337"

# Row 4
$ws.Range("G4").Value = [double]"1.046867666666667"
$ws.Range("H4").Value = [double]"3.140603"
$ws.Range("I4").Value = [double]"0.000687505225377314"
$ws.Range("J4").Value = [double]"0.000687505225377314"
$ws.Range("M4").Value = [double]"1.168795666666667"
$ws.Range("N4").Value = [double]"3.506387"
$ws.Range("O4").Value = [double]"0.3681476938764561"
$ws.Range("P4").Value = [double]"0.3681476938764561"
$ws.Range("Q4").Value = [double]"1.223574392373445"
$ws.Range("R4").Value = [double]"11.012169531361"
$ws.Range("S4").Value = [double]"0.0002531034632506713"
$ws.Range("T4").Value = [double]"0.0002531034632506713"

# Row 5
$ws.Range("H5").Value = [double]"4442.55542"
$ws.Range("I5").Value = [double]"0.9725138978974124"
$ws.Range("J5").Value = [double]"0.9725138978974125"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.2615913333333333"
$ws.Range("N5").Value = [double]"0.784774"
$ws.Range("O5").Value = [double]"0.08239613548481725"
$ws.Range("P5").Value = [double]"0.08239613548481727"
$ws.Range("Q5").Value = [double]"387.3779985750089"
$ws.Range("R5").Value = [double]"3486.40198717508"
$ws.Range("S5").Value = [double]"0.08013138689202293"
$ws.Range("T5").Value = [double]"0.08013138689202294"

# Row 6
$ws.Range("H6").Value = [double]"4442.55542"
$ws.Range("I6").Value = [double]"0.9725138978974124"
$ws.Range("J6").Value = [double]"0.9725138978974125"
$ws.Range("N6").Value = [double]"5.233242000000001"
$ws.Range("O6").Value = [double]"0.5494561706387266"
$ws.Range("P6").Value = [double]"0.5494561706387268"
$ws.Range("S6").Value = [double]"0.5343537622316538"
$ws.Range("T6").Value = [double]"0.534353762231654"

# Row 7
$ws.Range("H7").Value = [double]"4442.55542"
$ws.Range("I7").Value = [double]"0.9725138978974124"
$ws.Range("J7").Value = [double]"0.9725138978974125"
$ws.Range("M7").Value = [double]"1.168795666666667"
$ws.Range("N7").Value = [double]"3.506387"
$ws.Range("O7").Value = [double]"0.3681476938764561"
$ws.Range("P7").Value = [double]"0.3681476938764561"
$ws.Range("Q7").Value = [double]"1730.813174607505"
$ws.Range("R7").Value = [double]"15577.31857146754"
$ws.Range("S7").Value = [double]"0.3580287487737356"
$ws.Range("T7").Value = [double]"0.3580287487737357"

# Row 8
$ws.Range("G8").Value = [double]"40.80635833333333"
$ws.Range("H8").Value = [double]"122.419075"
$ws.Range("I8").Value = [double]"0.02679859687721029"
$ws.Range("J8").Value = [double]"0.0267985968772103"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"0.2615913333333333"
$ws.Range("N8").Value = [double]"0.784774"
$ws.Range("O8").Value = [double]"0.08239613548481725"
$ws.Range("P8").Value = [double]"0.08239613548481727"
$ws.Range("Q8").Value = [double]"10.67458968489444"
$ws.Range("R8").Value = [double]"96.07130716405"
$ws.Range("S8").Value = [double]"0.00220810081909762"
$ws.Range("T8").Value = [double]"0.002208100819097621"

# Row 9
$ws.Range("G9").Value = [double]"40.80635833333333"
$ws.Range("H9").Value = [double]"122.419075"
$ws.Range("I9").Value = [double]"0.02679859687721029"
$ws.Range("J9").Value = [double]"0.0267985968772103"
$ws.Range("N9").Value = [double]"5.233242000000001"
$ws.Range("O9").Value = [double]"0.5494561706387266"
$ws.Range("P9").Value = [double]"0.5494561706387268"
$ws.Range("Q9").Value = [double]"71.18318276568333"
$ws.Range("R9").Value = [double]"640.6486448911501"
$ws.Range("S9").Value = [double]"0.01472465441864291"
$ws.Range("T9").Value = [double]"0.01472465441864291"

# Row 10
$ws.Range("G10").Value = [double]"40.80635833333333"
$ws.Range("H10").Value = [double]"122.419075"
$ws.Range("I10").Value = [double]"0.02679859687721029"
$ws.Range("J10").Value = [double]"0.0267985968772103"
$ws.Range("M10").Value = [double]"1.168795666666667"
$ws.Range("N10").Value = [double]"3.506387"
$ws.Range("O10").Value = [double]"0.3681476938764561"
$ws.Range("P10").Value = [double]"0.3681476938764561"
$ws.Range("Q10").Value = [double]"47.69429479244722"
$ws.Range("R10").Value = [double]"429.248653132025"
$ws.Range("S10").Value = [double]"0.009865841639469767"
$ws.Range("T10").Value = [double]"0.009865841639469771"

